$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E contain numeric-looking / percentage text values that must
# remain plain text (matching the source inlineStr cells) rather than being
# auto-converted to numbers by Excel. Temporarily force Text format, write the
# values, then clear the format back so no stray style id is left behind.
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = "28.637.93"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").Value = "1.797.98"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "231.46"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").Value = "0.5895"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.2766"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "0.06798"
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("D10").Value = "23.22"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").Value = "0.07516"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "1.800.04"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "4.763"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("D14").Value = "0.6188"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "2.041.91"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "0.000009115"
$ws.Range("E16").Value = "  -8.02%  "
$ws.Range("D17").Value = "75.49"
$ws.Range("E17").Value = "  -4.78%  "
$ws.Range("D18").Value = "28.613.51"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "5.463"
$ws.Range("E19").Value = "  -6.58%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "210.48"
$ws.Range("E21").Value = "  -6.63%  "
$ws.Range("D22").Value = "11.49"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "6.811"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "153.52"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "7.861"
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("D27").Value = "0.1268"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("D28").Value = "16.42"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "1.423"
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("D30").Value = "0.06155"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").Value = "3.798"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "3.779"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "1.731"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").Value = "1.055"
$ws.Range("E35").Value = "  -6.22%  "
$ws.Range("D36").Value = "0.6412"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "2.500"
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("D38").Value = "2.714"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").Value = "6.536"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "0.01699"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("D41").Value = "1.147.39"
$ws.Range("D42").Value = "0.8866"
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("D43").Value = "1.007"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "100.09"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "1.944.94"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "60.25"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").Value = "0.00000000111"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").Value = "1.588"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.341"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05465"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "0.4480"
$ws.Range("E51").Value = "  -1.94%  "

$numRng.ClearFormats()
